$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 3).Value = 4.921302983646869
$ws.Cells.Item(2, 4).Value = 4.019281312647894
$ws.Cells.Item(2, 5).Value = 11.81046047891053
$ws.Cells.Item(2, 6).Value = 20.68889916068396
$ws.Cells.Item(2, 7).Value = 22.79351139556342
$ws.Cells.Item(2, 8).Value = 12.600286218835
$ws.Cells.Item(2, 11).Value = 11.08017004346835
$ws.Cells.Item(2, 13).Value = 14.23151592383905
$ws.Cells.Item(2, 14).Value = 16.5805041349021
$ws.Cells.Item(2, 15).Value = 18.40551535038954

$ws.Cells.Item(3, 3).Value = 4.74957993242162
$ws.Cells.Item(3, 4).Value = 3.95813619372516
$ws.Cells.Item(3, 5).Value = 11.70854087759793
$ws.Cells.Item(3, 6).Value = 20.65388348417508
$ws.Cells.Item(3, 7).Value = 22.6961987227487
$ws.Cells.Item(3, 8).Value = 12.63782602291546
$ws.Cells.Item(3, 11).Value = 10.52095607215808
$ws.Cells.Item(3, 13).Value = 13.94239640711534
$ws.Cells.Item(3, 14).Value = 16.6186984146878
$ws.Cells.Item(3, 15).Value = 18.44096003738967

$ws.Cells.Item(4, 3).Value = 4.642368593939826
$ws.Cells.Item(4, 4).Value = 3.919627257571541
$ws.Cells.Item(4, 5).Value = 11.65034441421413
$ws.Cells.Item(4, 6).Value = 20.63924142655172
$ws.Cells.Item(4, 7).Value = 22.64624179187773
$ws.Cells.Item(4, 8).Value = 12.66338828950767
$ws.Cells.Item(4, 11).Value = 10.15990018180379
$ws.Cells.Item(4, 13).Value = 13.76479395748087
$ws.Cells.Item(4, 14).Value = 16.64387955755739
$ws.Cells.Item(4, 15).Value = 18.46791433421766

$ws.Cells.Item(5, 3).Value = 4.598312170439432
$ws.Cells.Item(5, 4).Value = 3.903701856279045
$ws.Cells.Item(5, 5).Value = 11.6277553744019
$ws.Cells.Item(5, 6).Value = 20.63500272172733
$ws.Cells.Item(5, 7).Value = 22.62836227421622
$ws.Cells.Item(5, 8).Value = 12.67443583670337
$ws.Cells.Item(5, 11).Value = 10.00841213178984
$ws.Cells.Item(5, 13).Value = 13.69250215761628
$ws.Cells.Item(5, 14).Value = 16.65457671512079
$ws.Cells.Item(5, 15).Value = 18.48019946660192

$ws.Cells.Item(6, 3).Value = 4.590976902914472
$ws.Cells.Item(6, 4).Value = 3.901043703277796
$ws.Cells.Item(6, 5).Value = 11.62407317932476
$ws.Cells.Item(6, 6).Value = 20.63440332541773
$ws.Cells.Item(6, 7).Value = 22.62554344654709
$ws.Cells.Item(6, 8).Value = 12.676308333643
$ws.Cells.Item(6, 11).Value = 9.982997776274143
$ws.Cells.Item(6, 13).Value = 13.68050620804364
$ws.Cells.Item(6, 14).Value = 16.65637930194109
$ws.Cells.Item(6, 15).Value = 18.48231784716703

$ws.Cells.Item(7, 3).Value = 4.641775810382754
$ws.Cells.Item(7, 4).Value = 3.919413410626607
$ws.Cells.Item(7, 5).Value = 11.65003517914684
$ws.Cells.Item(7, 6).Value = 20.63917726177508
$ws.Cells.Item(7, 7).Value = 22.64599061051542
$ws.Cells.Item(7, 8).Value = 12.66353472837072
$ws.Cells.Item(7, 11).Value = 10.15787465455756
$ws.Cells.Item(7, 13).Value = 13.76381853244716
$ws.Cells.Item(7, 14).Value = 16.64402205836895
$ws.Cells.Item(7, 15).Value = 18.46807475404751

$ws.Cells.Item(8, 3).Value = 4.862506075962624
$ws.Cells.Item(8, 4).Value = 3.998404835094979
$ws.Cells.Item(8, 5).Value = 11.77442450113053
$ws.Cells.Item(8, 6).Value = 20.67540520640143
$ws.Cells.Item(8, 7).Value = 22.7579361981757
$ws.Cells.Item(8, 8).Value = 12.61270778015788
$ws.Cells.Item(8, 11).Value = 10.89107744523185
$ws.Cells.Item(8, 13).Value = 14.131905737837
$ws.Cells.Item(8, 14).Value = 16.59331493479851
$ws.Cells.Item(8, 15).Value = 18.41665681835864

$ws.Cells.Item(9, 3).Value = 5.27816186573907
$ws.Cells.Item(9, 4).Value = 4.145220795046396
$ws.Cells.Item(9, 5).Value = 12.05192029866666
$ws.Cells.Item(9, 6).Value = 20.80061098629886
$ws.Cells.Item(9, 7).Value = 23.05425240804665
$ws.Cells.Item(9, 8).Value = 12.53302165336108
$ws.Cells.Item(9, 11).Value = 12.18515189761282
$ws.Cells.Item(9, 13).Value = 14.84866919298516
$ws.Cells.Item(9, 14).Value = 16.50757260420576
$ws.Cells.Item(9, 15).Value = 18.35718674255263

$ws.Cells.Item(10, 3).Value = 5.569336587716885
$ws.Cells.Item(10, 4).Value = 4.247617588009789
$ws.Cells.Item(10, 5).Value = 12.27448084797895
$ws.Cells.Item(10, 6).Value = 20.92515269961422
$ws.Cells.Item(10, 7).Value = 23.31722958267576
$ws.Cells.Item(10, 8).Value = 12.4867237661105
$ws.Cells.Item(10, 11).Value = 13.04482813628736
$ws.Cells.Item(10, 13).Value = 15.36653580727707
$ws.Cells.Item(10, 14).Value = 16.45288557388182
$ws.Cells.Item(10, 15).Value = 18.33890992031236

$ws.Cells.Item(11, 3).Value = 5.698056683464816
$ws.Cells.Item(11, 4).Value = 4.29290704391937
$ws.Cells.Item(11, 5).Value = 12.37936740918842
$ws.Cells.Item(11, 6).Value = 20.9887377783956
$ws.Cells.Item(11, 7).Value = 23.44627827932862
$ws.Cells.Item(11, 8).Value = 12.46833371547136
$ws.Cells.Item(11, 11).Value = 13.41563973491201
$ws.Cells.Item(11, 13).Value = 15.59911497728084
$ws.Cells.Item(11, 14).Value = 16.42980272734409
$ws.Cells.Item(11, 15).Value = 18.33614534344226

$ws.Cells.Item(12, 3).Value = 5.746215687797622
$ws.Cells.Item(12, 4).Value = 4.309862956260027
$ws.Cells.Item(12, 5).Value = 12.41957181466246
$ws.Cells.Item(12, 6).Value = 21.01379749274995
$ws.Cells.Item(12, 7).Value = 23.49645936839601
$ws.Cells.Item(12, 8).Value = 12.46175492637425
$ws.Cells.Item(12, 11).Value = 13.55311267159467
$ws.Cells.Item(12, 13).Value = 15.68666886370023
$ws.Cells.Item(12, 14).Value = 16.42131929710293
$ws.Cells.Item(12, 15).Value = 18.33589801509326

$ws.Cells.Item(13, 3).Value = 5.735870471938187
$ws.Cells.Item(13, 4).Value = 4.306219966149253
$ws.Cells.Item(13, 5).Value = 12.41089206556737
$ws.Cells.Item(13, 6).Value = 21.00835704684131
$ws.Cells.Item(13, 7).Value = 23.48559426429982
$ws.Cells.Item(13, 8).Value = 12.46315464224214
$ws.Cells.Item(13, 11).Value = 13.52363684522929
$ws.Cells.Item(13, 13).Value = 15.66783697430825
$ws.Cells.Item(13, 14).Value = 16.42313490929946
$ws.Cells.Item(13, 15).Value = 18.33591570376948

$ws.Cells.Item(14, 3).Value = 5.702030718247642
$ws.Cells.Item(14, 4).Value = 4.294305956749275
$ws.Cells.Item(14, 5).Value = 12.38266552612622
$ws.Cells.Item(14, 6).Value = 20.99077988039105
$ws.Cells.Item(14, 7).Value = 23.45038063902339
$ws.Cells.Item(14, 8).Value = 12.46778475045204
$ws.Cells.Item(14, 11).Value = 13.42700890373686
$ws.Cells.Item(14, 13).Value = 15.60632896664877
$ws.Cells.Item(14, 14).Value = 16.42909963177117
$ws.Cells.Item(14, 15).Value = 18.33610896484305

$ws.Cells.Item(15, 3).Value = 5.681225468799435
$ws.Cells.Item(15, 4).Value = 4.286982743277949
$ws.Cells.Item(15, 5).Value = 12.36543813044648
$ws.Cells.Item(15, 6).Value = 20.98014068723964
$ws.Cells.Item(15, 7).Value = 23.42898100294332
$ws.Cells.Item(15, 8).Value = 12.47067100916273
$ws.Cells.Item(15, 11).Value = 13.36743709788063
$ws.Cells.Item(15, 13).Value = 15.56858347381792
$ws.Cells.Item(15, 14).Value = 16.43278672124557
$ws.Cells.Item(15, 15).Value = 18.33633150297511

$ws.Cells.Item(16, 3).Value = 5.560845072694863
$ws.Cells.Item(16, 4).Value = 4.244631090975625
$ws.Cells.Item(16, 5).Value = 12.26769636818562
$ws.Cells.Item(16, 6).Value = 20.92113546967013
$ws.Cells.Item(16, 7).Value = 23.30898203165706
$ws.Cells.Item(16, 8).Value = 12.48797943746038
$ws.Cells.Item(16, 11).Value = 13.02018414177494
$ws.Cells.Item(16, 13).Value = 15.35126824466379
$ws.Cells.Item(16, 14).Value = 16.45443015911152
$ws.Cells.Item(16, 15).Value = 18.33920240998219

$ws.Cells.Item(17, 3).Value = 5.486003156993373
$ws.Cells.Item(17, 4).Value = 4.218312766654994
$ws.Cells.Item(17, 5).Value = 12.20864148071719
$ws.Cells.Item(17, 6).Value = 20.88670253713009
$ws.Cells.Item(17, 7).Value = 23.23775073334424
$ws.Cells.Item(17, 8).Value = 12.49928247832246
$ws.Cells.Item(17, 11).Value = 12.80194153787989
$ws.Cells.Item(17, 13).Value = 15.21712001914969
$ws.Cells.Item(17, 14).Value = 16.46816696687183
$ws.Cells.Item(17, 15).Value = 18.34238630405152

$ws.Cells.Item(18, 3).Value = 5.442606790643204
$ws.Cells.Item(18, 4).Value = 4.203054095115915
$ws.Cells.Item(18, 5).Value = 12.1750198683752
$ws.Cells.Item(18, 6).Value = 20.86755097852473
$ws.Cells.Item(18, 7).Value = 23.19766964741035
$ws.Cells.Item(18, 8).Value = 12.50603504704999
$ws.Cells.Item(18, 11).Value = 12.67450864795659
$ws.Cells.Item(18, 13).Value = 15.13968393589567
$ws.Cells.Item(18, 14).Value = 16.47623695723711
$ws.Cells.Item(18, 15).Value = 18.34473990335316

$ws.Cells.Item(19, 3).Value = 5.427855068546394
$ws.Cells.Item(19, 4).Value = 4.197867244845656
$ws.Cells.Item(19, 5).Value = 12.16369658611259
$ws.Cells.Item(19, 6).Value = 20.86117923526217
$ws.Cells.Item(19, 7).Value = 23.18425281473846
$ws.Cells.Item(19, 8).Value = 12.50836448956746
$ws.Cells.Item(19, 11).Value = 12.63103588745692
$ws.Cells.Item(19, 13).Value = 15.11342048058683
$ws.Cells.Item(19, 14).Value = 16.47899834912915
$ws.Cells.Item(19, 15).Value = 18.34562644149051

$ws.Cells.Item(20, 3).Value = 5.494006715022285
$ws.Cells.Item(20, 4).Value = 4.221126990423624
$ws.Cells.Item(20, 5).Value = 12.21489251725463
$ws.Cells.Item(20, 6).Value = 20.89030046583771
$ws.Cells.Item(20, 7).Value = 23.24524167835061
$ws.Cells.Item(20, 8).Value = 12.49805322800771
$ws.Cells.Item(20, 11).Value = 12.82537130518994
$ws.Cells.Item(20, 13).Value = 15.23142970332026
$ws.Cells.Item(20, 14).Value = 16.46668717913882
$ws.Cells.Item(20, 15).Value = 18.34199330309444

$ws.Cells.Item(21, 3).Value = 5.711986485489391
$ws.Cells.Item(21, 4).Value = 4.297810726930003
$ws.Cells.Item(21, 5).Value = 12.39094345638071
$ws.Cells.Item(21, 6).Value = 20.99591621467959
$ws.Cells.Item(21, 7).Value = 23.46068844201947
$ws.Cells.Item(21, 8).Value = 12.46641431545749
$ws.Cells.Item(21, 11).Value = 13.45547102865434
$ws.Cells.Item(21, 13).Value = 15.6244101086135
$ws.Cells.Item(21, 14).Value = 16.42734066256229
$ws.Cells.Item(21, 15).Value = 18.33603049068762

$ws.Cells.Item(22, 3).Value = 5.85102393141529
$ws.Cells.Item(22, 4).Value = 4.346792274520242
$ws.Cells.Item(22, 5).Value = 12.50881842044491
$ws.Cells.Item(22, 6).Value = 21.07065438551432
$ws.Cells.Item(22, 7).Value = 23.60912789301464
$ws.Cells.Item(22, 8).Value = 12.44798168043299
$ws.Cells.Item(22, 11).Value = 13.85009975845549
$ws.Cells.Item(22, 13).Value = 15.87818487314298
$ws.Cells.Item(22, 14).Value = 16.40312636921344
$ws.Cells.Item(22, 15).Value = 18.33679412567219

$ws.Cells.Item(23, 3).Value = 5.777144563159795
$ws.Cells.Item(23, 4).Value = 4.320756517449249
$ws.Cells.Item(23, 5).Value = 12.44566140105676
$ws.Cells.Item(23, 6).Value = 21.0302480365032
$ws.Cells.Item(23, 7).Value = 23.52921909886533
$ws.Cells.Item(23, 8).Value = 12.45761376603367
$ws.Cells.Item(23, 11).Value = 13.64105953218002
$ws.Cells.Item(23, 13).Value = 15.74304800478574
$ws.Cells.Item(23, 14).Value = 16.41591282335155
$ws.Cells.Item(23, 15).Value = 18.3359597709335

$ws.Cells.Item(24, 3).Value = 5.490389451396347
$ws.Cells.Item(24, 4).Value = 4.219855077633369
$ws.Cells.Item(24, 5).Value = 12.21206539340119
$ws.Cells.Item(24, 6).Value = 20.88867183388621
$ws.Cells.Item(24, 7).Value = 23.24185230671673
$ws.Cells.Item(24, 8).Value = 12.49860818000667
$ws.Cells.Item(24, 11).Value = 12.81478481709168
$ws.Cells.Item(24, 13).Value = 15.22496126171022
$ws.Cells.Item(24, 14).Value = 16.46735565374131
$ws.Cells.Item(24, 15).Value = 18.34216934939203

$ws.Cells.Item(25, 3).Value = 5.167966922067204
$ws.Cells.Item(25, 4).Value = 4.106426130059
$ws.Cells.Item(25, 5).Value = 11.97343216735634
$ws.Cells.Item(25, 6).Value = 20.76098121002729
$ws.Cells.Item(25, 7).Value = 22.96601185893479
$ws.Cells.Item(25, 8).Value = 12.55243247238385
$ws.Cells.Item(25, 11).Value = 11.85086252460971
$ws.Cells.Item(25, 13).Value = 14.65590657925623
$ws.Cells.Item(25, 14).Value = 16.52930636710363
$ws.Cells.Item(25, 15).Value = 18.36882399714512

